$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Influencer")

# Update the header text for column N from "Chi phí Fb" to "Chi phí Facebook"
$ws.Range("N1").Value = "Chi phí Facebook"
